$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 576 (shifts existing rows 576:640 down to 577:641)
$ws.Rows.Item(576).Insert()

$ws.Cells.Item(576, 1).Value = 4
$ws.Cells.Item(576, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(576, 3).Value = "Los Lagos"
$ws.Cells.Item(576, 4).Value = 45194
$ws.Cells.Item(576, 5).Value = 10
$ws.Cells.Item(576, 6).Value = 100114013
$ws.Cells.Item(576, 7).Value = "Zanahoria"
$ws.Cells.Item(576, 8).Value = "Sin especificar"
$ws.Cells.Item(576, 9).Value = "Primera"
$ws.Cells.Item(576, 10).Value = 150
$ws.Cells.Item(576, 11).Value = 7500
$ws.Cells.Item(576, 12).Value = 7500
$ws.Cells.Item(576, 13).Value = 7500
$ws.Cells.Item(576, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(576, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(576, 16).Value = 375
$ws.Cells.Item(576, 17).Value = 20
$ws.Cells.Item(576, 18).Value = "Hortaliza"
